$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the year table with a new 2023 column (AA), carrying over the
# same per-row formatting that column Z (2022) already uses.
$ws.Range("Z4:Z16").Copy()
$ws.Range("AA4").PasteSpecial(-4122)

# New 2023 data.
$ws.Range("AA4").Value = 2023
$ws.Range("AA5").Value = 44.2
$ws.Range("AA6").Value = 50.4
$ws.Range("AA7").Value = 40.6
$ws.Range("AA8").Value = 57.2
$ws.Range("AA9").Value = 31
$ws.Range("AA10").Value = 49.7
$ws.Range("AA11").Value = 51
$ws.Range("AA12").Value = 29.4
$ws.Range("AA13").Value = 29.9
$ws.Range("AA14").Value = 56.3
$ws.Range("AA15").Value = 62.5
$ws.Range("AA16").Value = 34.9

# Reset the view: scroll back to the left edge and select A1 (the sheet
# had been scrolled to show the newly added column with it selected).
$ws.Range("A1").Select() | Out-Null
